$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4184.2563
$ws.Range("I64").Value = 4470.16
$ws.Range("J64").Value = 3673.7144
$ws.Range("K64").Value = 4470.16
$ws.Range("L64").Value = 3673.7144
$ws.Range("M64").Value = -4222.16
$ws.Range("N64").Value = -4169.7144

$ws.Range("H67").Value = 4184.2563
$ws.Range("I67").Value = 4470.16
$ws.Range("J67").Value = 3673.7144
$ws.Range("K67").Value = 4470.16
$ws.Range("L67").Value = 3673.7144
$ws.Range("M67").Value = -3612.16
$ws.Range("N67").Value = -5389.7144

$ws.Range("H74").Value = 6145.2
$ws.Range("I74").Value = 2725
$ws.Range("J74").Value = 8425.333
$ws.Range("K74").Value = 2725
$ws.Range("L74").Value = 8425.333
$ws.Range("M74").Value = -1789
$ws.Range("N74").Value = -10297.333

$ws.Range("H77").Value = 6145.2
$ws.Range("I77").Value = 2725
$ws.Range("J77").Value = 8425.333
$ws.Range("K77").Value = 13625
$ws.Range("L77").Value = 42126.665
$ws.Range("M77").Value = -8945
$ws.Range("N77").Value = -51486.665

$ws.Range("H80").Value = 525.38464
$ws.Range("I80").Value = 613.6667
$ws.Range("J80").Value = 498.9
$ws.Range("K80").Value = 1841.0001
$ws.Range("L80").Value = 1496.7
$ws.Range("M80").Value = -843.0001
$ws.Range("N80").Value = -3492.7

$ws.Range("H83").Value = 525.38464
$ws.Range("I83").Value = 613.6667
$ws.Range("J83").Value = 498.9
$ws.Range("K83").Value = 5523.0003
$ws.Range("L83").Value = 4490.099999999999
$ws.Range("M83").Value = -531.0002999999997
$ws.Range("N83").Value = -14474.1

$ws.Range("H103").Value = 813143
$ws.Range("I103").Value = 624
$ws.Range("J103").Value = 3250700
$ws.Range("K103").Value = 1872
$ws.Range("L103").Value = 9752100
$ws.Range("M103").Value = -1286
$ws.Range("N103").Value = -9753272

$ws.Range("H112").Value = 11112020
$ws.Range("J112").Value = 13889847
$ws.Range("L112").Value = 41669541
$ws.Range("N112").Value = -41671757

$ws.Range("H137").Value = 1644.975
$ws.Range("I137").Value = 1181.3334
$ws.Range("J137").Value = 2024.3182
$ws.Range("K137").Value = 3544.0002
$ws.Range("L137").Value = 6072.9546
$ws.Range("M137").Value = -994.0002
$ws.Range("N137").Value = -11172.9546

$ws.Range("H138").Value = 2928.9402
$ws.Range("I138").Value = 1287.9395
$ws.Range("J138").Value = 4521.6763
$ws.Range("K138").Value = 3863.8185
$ws.Range("L138").Value = 13565.0289
$ws.Range("M138").Value = 1276.1815
$ws.Range("N138").Value = -23845.0289

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 125002216
$ws.Range("I88").Value = 2526
$ws.Range("J88").Value = 250001900
$ws.Range("K88").Value = 2526
$ws.Range("L88").Value = 250001900
$ws.Range("M88").Value = -2120
$ws.Range("N88").Value = -250002712

$ws.Range("H91").Value = 125002216
$ws.Range("I91").Value = 2526
$ws.Range("J91").Value = 250001900
$ws.Range("K91").Value = 2526
$ws.Range("L91").Value = 250001900
$ws.Range("M91").Value = -1122
$ws.Range("N91").Value = -250004708

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 9983.143
$ws.Range("I54").Value = 6976.4
$ws.Range("J54").Value = 17500
$ws.Range("K54").Value = 6976.4
$ws.Range("L54").Value = 17500
$ws.Range("M54").Value = -6492.4
$ws.Range("N54").Value = -18468

$ws.Range("H86").Value = 2779.6
$ws.Range("I86").Value = 2688.4443
$ws.Range("K86").Value = 2688.4443
$ws.Range("M86").Value = -1565.4443

$ws.Range("H89").Value = 2779.6
$ws.Range("I89").Value = 2688.4443
$ws.Range("K89").Value = 13442.2215
$ws.Range("M89").Value = -7826.2215

$ws.Range("H94").Value = 1503.8928
$ws.Range("I94").Value = 1257.6666
$ws.Range("K94").Value = 1257.6666
$ws.Range("M94").Value = -806.6666

$ws.Range("H107").Value = 1001.61536
$ws.Range("I107").Value = 1026.25
$ws.Range("J107").Value = 962.2
$ws.Range("K107").Value = 1026.25
$ws.Range("L107").Value = 962.2
$ws.Range("M107").Value = 893.75
$ws.Range("N107").Value = -4802.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3726.81
$ws.Range("I31").Value = 1838.2273
$ws.Range("J31").Value = 4455.737
$ws.Range("K31").Value = 1838.2273
$ws.Range("L31").Value = 4455.737
$ws.Range("M31").Value = -1543.2273
$ws.Range("N31").Value = -5045.737

$ws.Range("H34").Value = 3726.81
$ws.Range("I34").Value = 1838.2273
$ws.Range("J34").Value = 4455.737
$ws.Range("K34").Value = 1838.2273
$ws.Range("L34").Value = 4455.737
$ws.Range("M34").Value = -1636.2273
$ws.Range("N34").Value = -4859.737

$ws.Range("H48").Value = 24999
$ws.Range("J48").Value = 24999
$ws.Range("L48").Value = 24999
$ws.Range("N48").Value = -25951

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4850.5
$ws.Range("J55").Value = 5333.8887
$ws.Range("L55").Value = 16001.6661
$ws.Range("N55").Value = -16355.6661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6062.9023
$ws.Range("I70").Value = 6130.2905
$ws.Range("K70").Value = 6130.2905
$ws.Range("M70").Value = -5860.2905

$ws.Range("H73").Value = 6062.9023
$ws.Range("I73").Value = 6130.2905
$ws.Range("K73").Value = 6130.2905
$ws.Range("M73").Value = -5194.2905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 9496.25
$ws.Range("I58").Value = 3992.5
$ws.Range("J58").Value = 15000
$ws.Range("K58").Value = 3992.5
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = -3732.5
$ws.Range("N58").Value = -15520

$ws.Range("H100").Value = 2342
$ws.Range("J100").Value = 2401.3333
$ws.Range("L100").Value = 2401.3333
$ws.Range("N100").Value = -3483.3333

$ws.Range("H122").Value = 5095435
$ws.Range("I122").Value = 7147535
$ws.Range("K122").Value = 21442605
$ws.Range("M122").Value = -21440155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 7471
$ws.Range("J54").Value = 7471
$ws.Range("L54").Value = 7471
$ws.Range("N54").Value = -8511

$ws.Range("H81").Value = 2675.3333
$ws.Range("J81").Value = 2675.3333
$ws.Range("L81").Value = 5350.6666
$ws.Range("N81").Value = -7472.6666

$ws.Range("H84").Value = 2675.3333
$ws.Range("J84").Value = 2675.3333
$ws.Range("L84").Value = 26753.333
$ws.Range("N84").Value = -37361.333
